$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 51
$ws.Range("C2").Value = 56
$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 46
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 15
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 10
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 30
$ws.Range("B11").Value = 47
$ws.Range("C11").Value = 56
$ws.Range("B17").Value = 28
$ws.Range("C17").Value = 16
$ws.Range("B18").Value = 13
$ws.Range("C18").Value = 22
$ws.Range("B20").Value = 27
$ws.Range("C20").Value = 24
$ws.Range("B22").Value = 31
$ws.Range("C22").Value = 24
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 25
$ws.Range("B27").Value = 21
$ws.Range("C27").Value = 27
$ws.Range("B28").Value = 41
$ws.Range("C28").Value = 43
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = 37
$ws.Range("B34").Value = 15
$ws.Range("C34").Value = 28
$ws.Range("B36").Value = 34
$ws.Range("C36").Value = 30
$ws.Range("B37").Value = 17
$ws.Range("C37").Value = 33
$ws.Range("B40").Value = 8
$ws.Range("C40").Value = 27
$ws.Range("B41").Value = 34
$ws.Range("C41").Value = 40
$ws.Range("B43").Value = 36
$ws.Range("C43").Value = 47
$ws.Range("B45").Value = 32
$ws.Range("C45").Value = 19
$ws.Range("B46").Value = 9
$ws.Range("C46").Value = 30
$ws.Range("B47").Value = 15
$ws.Range("C47").Value = 37
$ws.Range("B48").Value = 39
$ws.Range("C48").Value = 48
$ws.Range("B49").Value = 50
$ws.Range("C49").Value = 58
$ws.Range("B56").Value = 10
$ws.Range("C56").Value = 24
$ws.Range("B57").Value = 46
$ws.Range("C57").Value = 60
$ws.Range("B58").Value = 51
$ws.Range("C58").Value = 57
$ws.Range("B61").Value = 22
$ws.Range("C61").Value = 22
$ws.Range("B64").Value = 19
$ws.Range("C64").Value = 33
$ws.Range("B74").Value = 7
$ws.Range("C74").Value = 27
$ws.Range("B75").Value = 16
$ws.Range("C75").Value = 29
$ws.Range("B77").Value = 48
$ws.Range("C77").Value = 33
$ws.Range("B79").Value = 7
$ws.Range("C79").Value = 15
$ws.Range("B80").Value = 48
$ws.Range("C80").Value = 52
$ws.Range("B81").Value = 13
$ws.Range("C81").Value = 20
$ws.Range("B82").Value = 10
$ws.Range("C82").Value = 17
$ws.Range("B86").Value = 36
$ws.Range("C86").Value = 32
$ws.Range("B87").Value = 39
$ws.Range("C87").Value = 55
$ws.Range("B88").Value = 49
$ws.Range("C88").Value = 50
$ws.Range("B89").Value = 33
$ws.Range("C89").Value = 38
$ws.Range("B90").Value = 32
$ws.Range("C90").Value = 31
$ws.Range("B91").Value = 25
$ws.Range("C91").Value = 38
$ws.Range("B92").Value = 42
$ws.Range("C92").Value = 39
$ws.Range("B97").Value = 28
$ws.Range("C97").Value = 43
$ws.Range("B98").Value = 9
$ws.Range("C98").Value = 29
$ws.Range("B99").Value = 28
$ws.Range("C99").Value = 12
$ws.Range("B102").Value = 35
$ws.Range("C102").Value = 31
$ws.Range("B103").Value = 42
$ws.Range("C103").Value = 51
$ws.Range("B106").Value = 31
$ws.Range("C106").Value = 20
$ws.Range("B107").Value = 13
$ws.Range("C107").Value = 24
$ws.Range("B110").Value = 29
$ws.Range("C110").Value = 35
$ws.Range("B111").Value = 48
$ws.Range("C111").Value = 51
$ws.Range("B112").Value = 54
$ws.Range("C112").Value = 58
$ws.Range("B113").Value = 14
$ws.Range("C113").Value = 24
$ws.Range("B114").Value = 11
$ws.Range("C114").Value = 19
$ws.Range("B115").Value = 12
$ws.Range("C115").Value = 26
$ws.Range("B118").Value = 11
$ws.Range("C118").Value = 24
$ws.Range("B119").Value = 14
$ws.Range("C119").Value = 23
$ws.Range("B122").Value = 52
$ws.Range("C122").Value = 56
$ws.Range("B123").Value = 13
$ws.Range("C123").Value = 34
$ws.Range("B124").Value = 17
$ws.Range("C124").Value = 29
$ws.Range("B125").Value = 47
$ws.Range("C125").Value = 52
$ws.Range("B127").Value = 45
$ws.Range("C127").Value = 47
$ws.Range("B129").Value = 15
$ws.Range("C129").Value = 35
$ws.Range("B131").Value = 9
$ws.Range("C131").Value = 27
$ws.Range("B133").Value = 28
$ws.Range("C133").Value = 21
$ws.Range("B134").Value = 29
$ws.Range("C134").Value = 33
$ws.Range("B135").Value = 27
$ws.Range("C135").Value = 21
